# Update profit.py after running on 2025-09-01
#
# - Sheet1: append the new daily profit row for 09/01/2025.
# - Sheet2: refresh the single summary row to the latest date + ratios.
#
# Dates are written as literal text (matching the existing "MM/DD/2025"
# string cells already in the sheet), not as Excel date serials, so the
# cell is pre-formatted as Text before the value is assigned and the
# formatting is cleared again afterwards to avoid leaving a stray
# number-format behind on the cell.

$wb = $excel.ActiveWorkbook

# --- Sheet1: append new daily profit row ---
$ws1 = $wb.Worksheets.Item("Sheet1")

$newRow = 15
$dateCell1 = $ws1.Cells.Item($newRow, 1)
$dateCell1.NumberFormat = "@"
$dateCell1.Value = "09/01/2025"
$dateCell1.ClearFormats()

$ws1.Cells.Item($newRow, 2).Value = 11390.83

# --- Sheet2: update the summary row with the latest date + ratios ---
$ws2 = $wb.Worksheets.Item("Sheet2")

$dateCell2 = $ws2.Cells.Item(1, 1)
$dateCell2.NumberFormat = "@"
$dateCell2.Value = "09/01/2025"
$dateCell2.ClearFormats()

$ws2.Cells.Item(1, 2).Value = 0.1138573803187714
$ws2.Cells.Item(1, 3).Value = 0.8861426196812286
